$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 3571.6667
$ws.Range("I6").Value = 3571.6667
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 10715.0001
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -10603.0001
$ws.Range("N6").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1608.1666
$ws.Range("I40").Value = 1622.1111
$ws.Range("J40").Value = 1566.3334
$ws.Range("K40").Value = 1622.1111
$ws.Range("L40").Value = 1566.3334
$ws.Range("M40").Value = -1447.1111
$ws.Range("N40").Value = -1916.3334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 96.25
$ws.Range("I58").Value = 96.25
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 288.75
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -138.75
$ws.Range("N58").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 14603.6
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 14603.6
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 14603.6
$ws.Range("N113").Value = -21111.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("M135").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 2000
$ws.Range("I11").Value = 2000
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 2000
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -1856

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5222
$ws.Range("I45").Value = 5222
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 5222
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -4845

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3224.5
$ws.Range("I61").Value = 2966
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 2966
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -2754
$ws.Range("N61").Value = -4424

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1336666.4
$ws.Range("I74").Value = 1336666.4
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1336666.4
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -1335792.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1336666.4
$ws.Range("I77").Value = 1336666.4
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 6683332
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -6678964

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3224.5
$ws.Range("I136").Value = 2966
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 8898
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -6348
$ws.Range("N136").Value = -17100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 1208.375
$ws.Range("I12").Value = 845.2857
$ws.Range("J12").Value = 3750
$ws.Range("K12").Value = 845.2857
$ws.Range("L12").Value = 3750
$ws.Range("M12").Value = -677.2857

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5779.4
$ws.Range("I134").Value = 5779.4
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 17338.2
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -14803.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 700
$ws.Range("I4").Value = 2000
$ws.Range("J4").Value = 440
$ws.Range("K4").Value = 2000
$ws.Range("L4").Value = 440
$ws.Range("M4").Value = -1888
$ws.Range("N4").Value = -664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 9999
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 9999
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 9999
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -10339

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 57.5
$ws.Range("I22").Value = 57.5
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 57.5
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 292.5
$ws.Range("N22").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2000
$ws.Range("I31").Value = 2000
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 2000
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1705

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2000
$ws.Range("I34").Value = 2000
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 2000
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1798

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4868.2856
$ws.Range("I132").Value = 3651.9092
$ws.Range("J132").Value = 9328.333000000001
$ws.Range("K132").Value = 10955.7276
$ws.Range("L132").Value = 27984.999
$ws.Range("M132").Value = -8425.7276

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1211
$ws.Range("I134").Value = 1211
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 3633
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -1098
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 145000
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 145000
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 145000
$ws.Range("N135").Value = -155140

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 338.45456
$ws.Range("I8").Value = 338.45456
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 1015.36368
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -876.36368

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 174.1
$ws.Range("I10").Value = 27.235294
$ws.Range("J10").Value = 1006.3333
$ws.Range("K10").Value = 81.705882
$ws.Range("L10").Value = 3018.9999
$ws.Range("M10").Value = 57.294118
$ws.Range("N10").Value = -3296.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 3832
$ws.Range("I46").Value = 3
$ws.Range("J46").Value = 4597.8
$ws.Range("K46").Value = 9
$ws.Range("L46").Value = 13793.4
$ws.Range("M46").Value = 82
$ws.Range("N46").Value = -13975.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1336.5714
$ws.Range("I129").Value = 320
$ws.Range("J129").Value = 2099
$ws.Range("K129").Value = 960
$ws.Range("L129").Value = 6297
$ws.Range("M129").Value = 4040
$ws.Range("N129").Value = -16297

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 2348.75
$ws.Range("I137").Value = 1000
$ws.Range("J137").Value = 2798.3333
$ws.Range("K137").Value = 3000
$ws.Range("L137").Value = 8394.999899999999
$ws.Range("M137").Value = 2100
$ws.Range("N137").Value = -18594.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 77.42308
$ws.Range("I2").Value = 58.61111
$ws.Range("J2").Value = 119.75
$ws.Range("K2").Value = 58.61111
$ws.Range("L2").Value = 119.75
$ws.Range("M2").Value = 54.38889

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4734.8667
$ws.Range("I122").Value = 3844.4443
$ws.Range("J122").Value = 6070.5
$ws.Range("K122").Value = 11533.3329
$ws.Range("L122").Value = 18211.5
$ws.Range("M122").Value = -9083.332900000001
$ws.Range("N122").Value = -23111.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 357.76923
$ws.Range("I46").Value = 323.27274
$ws.Range("J46").Value = 547.5
$ws.Range("K46").Value = 323.27274
$ws.Range("L46").Value = 547.5
$ws.Range("M46").Value = -135.27274
$ws.Range("N46").Value = -923.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H105").Value = 29750
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 29750
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 29750
$ws.Range("N105").Value = -36738

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3461.4614
$ws.Range("I122").Value = 3180
$ws.Range("J122").Value = 3637.375
$ws.Range("K122").Value = 9540
$ws.Range("L122").Value = 10912.125
$ws.Range("M122").Value = -7090
$ws.Range("N122").Value = -15812.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 37455.5
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 37455.5
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 37455.5
$ws.Range("N103").Value = -39799.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1979.5
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 1979.5
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 5938.5
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -9778.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 439.2857
$ws.Range("I113").Value = 456.66666
$ws.Range("J113").Value = 426.25
$ws.Range("K113").Value = 1369.99998
$ws.Range("L113").Value = 1278.75
$ws.Range("M113").Value = 800.0000199999999
$ws.Range("N113").Value = -5618.75
